$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the input values (dependent formulas recalc automatically)
$ws.Range("B4").Value = 0.05
$ws.Range("B5").Value = 1

# Move the active selection to B4
$ws.Range("B4").Select()
